$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Fixed surplus number" - the surcharge rows that were still using the old
# "SS" (stainless steel) surcharge multiplier of 1.0565 are corrected to 1,
# matching the rest of the "Others" rows (K41/K42 legend: SS @ 1.0565 vs
# Others @ 1.0).
$surchargeCells = @("K16", "K17", "K20", "K23", "K26", "K29", "K32", "K35", "K36", "K40")
foreach ($cellRef in $surchargeCells) {
    $ws.Range($cellRef).Value = 1
}

# Reflect where the editor was last working in the sheet.
$ws.Range("K41").Select()
